# Commit: "Update posts.xlsx after post"
#
# The diff shows the data row for post "「自信を持つには？」" (row 226) being
# removed entirely, with every following row shifting up by one
# (dimension goes from A1:C357 to A1:C356). This is a straightforward
# whole-row deletion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A226").EntireRow.Delete()
